$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update the "publish date" header cells for the latest period column (G/H)
$ws.Range("G9").Value = "1402-04-14 (9)"
$ws.Range("H9").Value = "1402-04-28 (3)"

# Updated income-statement figures for the latest period (columns G/H)
$ws.Range("H12").Value = -411200
$ws.Range("H13").Value = 195826

$ws.Range("G14").Value = -115455
$ws.Range("H14").Value = -125219

$ws.Range("H16").Value = 16784

$ws.Range("G17").Value = 57719
$ws.Range("H17").Value = 87391

$ws.Range("H19").Value = -1359

$ws.Range("G20").Value = 71884
$ws.Range("H20").Value = 86032

$ws.Range("G22").Value = 71884
$ws.Range("H22").Value = 86032

$ws.Range("G24").Value = 71884
$ws.Range("H24").Value = 86032
